$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.659.67'
$ws.Range('E2').Value = '  -2.49%  '
$ws.Range('D3').Value = '2.292.00'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '546.77'
$ws.Range('E5').Value = '  -1.30%  '
$ws.Range('D6').Value = '130.53'
$ws.Range('E6').Value = '  -4.98%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '0.571'
$ws.Range('E8').Value = '  -3.04%  '
$ws.Range('D9').Value = '2.290.69'
$ws.Range('E10').Value = '  -3.71%  '
$ws.Range('D11').Value = '5.54'
$ws.Range('E11').Value = '  -3.24%  '
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('E13').Value = '  -5.20%  '
$ws.Range('D14').Value = '23.76'
$ws.Range('E14').Value = '  -4.69%  '
$ws.Range('D15').Value = '2.698.38'
$ws.Range('E15').Value = '  -5.35%  '
$ws.Range('D16').Value = '58.609.77'
$ws.Range('E16').Value = '  -2.43%  '
$ws.Range('E17').Value = '  -3.61%  '
$ws.Range('D18').Value = '2.229.88'
$ws.Range('E18').Value = '  -7.64%  '
$ws.Range('E19').Value = '  -5.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.30'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -4.45%  '
$ws.Range('D21').Value = '314.16'
$ws.Range('E21').Value = '  -4.08%  '
$ws.Range('E22').Value = '  -4.14%  '
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('E24').Value = '  -3.30%  '
$ws.Range('E25').Value = '  -3.76%  '
$ws.Range('D27').Value = '8.07'
$ws.Range('E27').Value = '  -7.30%  '
$ws.Range('D28').Value = '1.32'
$ws.Range('E28').Value = '  -6.63%  '
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('D30').Value = '169.28'
$ws.Range('E30').Value = '  -0.57%  '
$ws.Range('E31').Value = '  -6.38%  '
$ws.Range('D32').Value = '5.78'
$ws.Range('E32').Value = '  -5.39%  '
$ws.Range('E33').Value = '  +0.28%  '
$ws.Range('E34').Value = '  -5.23%  '
$ws.Range('D36').Value = '17.78'
$ws.Range('E36').Value = '  -4.14%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('D38').Value = '1.24'
$ws.Range('E38').Value = '  -6.42%  '
$ws.Range('D39').Value = '3.95'
$ws.Range('E39').Value = '  -6.46%  '
$ws.Range('D40').Value = '37.93'
$ws.Range('E40').Value = '  -2.40%  '
$ws.Range('E41').Value = '  -5.90%  '
$ws.Range('D42').Value = '299.31'
$ws.Range('E42').Value = '  -8.40%  '
$ws.Range('D43').Value = '140.12'
$ws.Range('E43').Value = '  -2.77%  '
$ws.Range('D44').Value = '3.44'
$ws.Range('E44').Value = '  -6.05%  '
$ws.Range('D45').Value = '0.0952'
$ws.Range('E45').Value = '  -1.41%  '
$ws.Range('E46').Value = '  -3.46%  '
$ws.Range('D47').Value = '0.557'
$ws.Range('E47').Value = '  -3.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.50'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -7.61%  '
$ws.Range('E49').Value = '  -3.86%  '
$ws.Range('E50').Value = '  -4.97%  '
$ws.Range('D51').Value = '11.01'
$ws.Range('E51').Value = '  -0.27%  '
